# StrategyExample.xlsx update: "data updated & ph added"
#
# - Fix two ticker labels that used a slash (which Excel/shared-strings had
#   stored as "NFT50SRTDURDEBTDYNP/B" and ".../E") to use an underscore
#   instead ("NFT50SRTDURDEBTDYNP_B" / "..._E"), and re-save so those
#   corrected strings end up after the already-present "NFTEQTSAV" entry
#   in the shared-string table (matches the row order A38:A40 -> NFTEQTSAV
#   row stays last/unchanged, the two corrected tickers take its old slot).
# - Slightly widen column A to fit the (now longer) label text.
# - Leave the active selection on D31 (scrolled so row 16 is near the top
#   of the view) instead of the old E37 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the ticker names in column A (rows 38-39): slash -> underscore.
# Row 40 already holds "NFTEQTSAV" and is unchanged.
$ws.Range("A38").Value = "NFT50SRTDURDEBTDYNP_B"
$ws.Range("A39").Value = "NFT50SRTDURDEBTDYNP_E"
$ws.Range("A40").Value = "NFTEQTSAV"

# --- Widen column A slightly so the longer text keeps its "best fit" look.
$ws.Columns("A").ColumnWidth = 23

# --- Update the view/selection: scroll to row 16 and select D31.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("D31").Select()
